$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.154684662818909
$ws.Range("B1").Value = 1.870872735977173
$ws.Range("D1").Value = 2.219284296035767
$ws.Range("E1").Value = 1.085615634918213
